# Add the "name:string" field to the devise User scaffold command.
# "rails g devise User"  ->  "rails g devise User name:string"

$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "rails g devise User",  # FindText
    $true,                  # MatchCase
    $true,                  # MatchWholeWord
    $false,                 # MatchWildcards
    $false,                 # MatchSoundsLike
    $false,                 # MatchAllWordForms
    $true,                  # Forward
    1,                      # Wrap (wdFindContinue)
    $false,                 # Format
    "rails g devise User name:string",  # ReplaceWith
    2                       # Replace (wdReplaceAll)
)
